$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 25000
$ws.Range("I2").Value = 20000
$ws.Range("P2").Value = 10000
$ws.Range("S2").Value = 73158
$ws.Range("T2").Value = 5416.988499999999
$ws.Range("U2").Value = 55000
$ws.Range("S3").Value = 69971
$ws.Range("T3").Value = 5210.142
$ws.Range("S4").Value = 67233
$ws.Range("T4").Value = 5131.699999999999
$ws.Range("S5").Value = 66774
$ws.Range("T5").Value = 5049.197999999999
$ws.Range("S6").Value = 68427
$ws.Range("T6").Value = 5084.9925
$ws.Range("S7").Value = 72833
$ws.Range("T7").Value = 5380.1055
$ws.Range("S8").Value = 72858
$ws.Range("T8").Value = 6413.1515
$ws.Range("S9").Value = 84685
$ws.Range("T9").Value = 7207.8405
$ws.Range("P10").Value = 10000
$ws.Range("S10").Value = 102552
$ws.Range("T10").Value = 9489.2
$ws.Range("U10").Value = 55000
$ws.Range("S11").Value = 112257
$ws.Range("T11").Value = 14318.1885
$ws.Range("S12").Value = 116157
$ws.Range("T12").Value = 15340.192
$ws.Range("S13").Value = 115181
$ws.Range("T13").Value = 15589.8715
$ws.Range("S14").Value = 118438
$ws.Range("T14").Value = 15359.3685
$ws.Range("S15").Value = 118968
$ws.Range("T15").Value = 15511.44
$ws.Range("S16").Value = 121806
$ws.Range("T16").Value = 15666.4375
$ws.Range("S17").Value = 114966
$ws.Range("T17").Value = 16125.417
$ws.Range("S18").Value = 15556
$ws.Range("T18").Value = 16406.873
$ws.Range("V18").Value = 3891.746977777778
$ws.Range("W18").Value = -98.93470541392529
$ws.Range("T19").Value = 16163.1645
$ws.Range("T20").Value = 15659.7175
$ws.Range("T21").Value = 13764.5865
$ws.Range("T22").Value = 11506.208
$ws.Range("T23").Value = 8102.793999999999
$ws.Range("T24").Value = 6308.112999999999
$ws.Range("T25").Value = 5711.482
